$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = $false
$ws.Range("F8").Value = $true
$ws.Range("F9").Value = $false

$ws.Range("A11").Value = "Current"
$ws.Range("B11").Value = $false
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = $true
$ws.Range("E11").Value = $true
$ws.Range("F11").Value = $false
$ws.Range("G11").Value = $true
$ws.Range("H11").Value = $true
$ws.Range("I11").Value = $true

$ws.Range("B9").Select()
